# edit.ps1
# Applies the "New crime data collected" update to the NYPD CompStat weekly
# workbook: advances the report week/volume text, and refreshes all of the
# weekly/28-day/YTD crime statistic counts and percentage changes for rows
# 14-30 (and the associated text/placeholder cells in rows 31 and 33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header text: bump the report volume/week number and the reporting
#    week's start/end dates.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/1/2024  Through  4/7/2024"

# ---------------------------------------------------------------------
# 2. Crime-complaint statistics table (rows 14-30): update the raw counts
#    and recalculated percentage changes for Week-to-Date, 28-Day,
#    Year-to-Date and the historical comparison columns.
# ---------------------------------------------------------------------

# Row 14 - Murder
$ws.Range("J14").Value = 6
$ws.Range("K14").Value = -83.333333333333

# Row 15 - Rape
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = -37.5
$ws.Range("I15").Value = 16
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = 23.076923076923
$ws.Range("L15").Value = 166.666666666667
$ws.Range("M15").Value = 45.454545454545
$ws.Range("N15").Value = -42.857142857142

# Row 16 - Robbery
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 32
$ws.Range("G16").Value = 26
$ws.Range("H16").Value = 23.076923076923
$ws.Range("I16").Value = 83
$ws.Range("J16").Value = 97
$ws.Range("K16").Value = -14.432989690721
$ws.Range("L16").Value = 31.746031746031
$ws.Range("M16").Value = -11.702127659574
$ws.Range("N16").Value = -77.445652173913

# Row 17 - Fel. Assault
$ws.Range("D17").Value = 23
$ws.Range("E17").Value = -34.782608695652
$ws.Range("F17").Value = 72
$ws.Range("G17").Value = 76
$ws.Range("H17").Value = -5.263157894736
$ws.Range("I17").Value = 262
$ws.Range("J17").Value = 235
$ws.Range("K17").Value = 11.489361702127
$ws.Range("L17").Value = 35.751295336787
$ws.Range("M17").Value = 92.647058823529
$ws.Range("N17").Value = -3.321033210332

# Row 18 - Burglary
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -28.571428571428
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 34
$ws.Range("H18").Value = -38.235294117647
$ws.Range("I18").Value = 69
$ws.Range("J18").Value = 103
$ws.Range("K18").Value = -33.009708737864
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -50.359712230215
$ws.Range("N18").Value = -92.123287671232

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 26
$ws.Range("D19").Value = 38
$ws.Range("E19").Value = -31.578947368421
$ws.Range("F19").Value = 116
$ws.Range("G19").Value = 120
$ws.Range("H19").Value = -3.333333333333
$ws.Range("I19").Value = 411
$ws.Range("J19").Value = 408
$ws.Range("K19").Value = 0.735294117647
$ws.Range("L19").Value = 8.730158730158
$ws.Range("M19").Value = 72.689075630252
$ws.Range("N19").Value = 3.007518796992

# Row 20 - G.L.A.
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -14.285714285714
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = -22.222222222222
$ws.Range("I20").Value = 63
$ws.Range("J20").Value = 96
$ws.Range("K20").Value = -34.375
$ws.Range("L20").Value = -36.363636363636
$ws.Range("M20").Value = -13.698630136986
$ws.Range("N20").Value = -95.112490302560

# Row 21 - TOTAL
$ws.Range("C21").Value = 61
$ws.Range("D21").Value = 87
$ws.Range("E21").Value = -29.885057471264
$ws.Range("F21").Value = 267
$ws.Range("G21").Value = 293
$ws.Range("H21").Value = -8.873720136518
$ws.Range("I21").Value = 905
$ws.Range("J21").Value = 958
$ws.Range("K21").Value = -5.532359081419
$ws.Range("L21").Value = 11.866501854140
$ws.Range("M21").Value = 30.591630591630
$ws.Range("N21").Value = -72.042014210688

# Row 23 - Housing
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = -50
$ws.Range("I23").Value = 26
$ws.Range("J23").Value = 45
$ws.Range("K23").Value = -42.222222222222
$ws.Range("L23").Value = -10.344827586206
$ws.Range("M23").Value = 85.714285714285

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 82
$ws.Range("D24").Value = 77
$ws.Range("E24").Value = 6.493506493506
$ws.Range("F24").Value = 354
$ws.Range("G24").Value = 303
$ws.Range("H24").Value = 16.831683168316
$ws.Range("I24").Value = 1191
$ws.Range("J24").Value = 1052
$ws.Range("K24").Value = 13.212927756654
$ws.Range("L24").Value = 35.803876852907
$ws.Range("M24").Value = 40.117647058823

# Row 25 - Retail Theft
$ws.Range("C25").Value = 43
$ws.Range("D25").Value = 42
$ws.Range("E25").Value = 2.380952380952
$ws.Range("F25").Value = 202
$ws.Range("H25").Value = 29.487179487179
$ws.Range("I25").Value = 675
$ws.Range("J25").Value = 518
$ws.Range("K25").Value = 30.308880308880
$ws.Range("L25").Value = 117.041800643087

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 28
$ws.Range("D26").Value = 44
$ws.Range("E26").Value = -36.363636363636
$ws.Range("F26").Value = 135
$ws.Range("G26").Value = 149
$ws.Range("H26").Value = -9.395973154362
$ws.Range("I26").Value = 443
$ws.Range("J26").Value = 454
$ws.Range("K26").Value = -2.422907488986
$ws.Range("L26").Value = 11.027568922305
$ws.Range("M26").Value = -6.144067796610

# Row 27 - UCR Rape*
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -40
$ws.Range("I27").Value = 28
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = 47.368421052631
$ws.Range("L27").Value = 64.705882352941

# Row 28 - Other Sex Crimes
$ws.Range("F28").Value = 13
$ws.Range("G28").Value = 19
$ws.Range("H28").Value = -31.578947368421
$ws.Range("I28").Value = 44
$ws.Range("J28").Value = 50
$ws.Range("K28").Value = -12
$ws.Range("L28").Value = -2.222222222222

# Row 29 - Shooting Vic.
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -66.666666666666
$ws.Range("J29").Value = 5
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = -58.333333333333

# Row 30 - Shooting Inc.
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = -66.666666666666
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = -50

# ---------------------------------------------------------------------
# 3. Cells that switch from a numeric value to the "no data" placeholder
#    text (shared-string "0" / "***.*"). A plain .Value assignment of a
#    numeric-looking string gets stored back as a number, so instead we
#    copy an existing cell that already has the desired text + style
#    onto the target cell - this preserves both the "s" (shared string)
#    cell type and the original style index (no new styles created).
# ---------------------------------------------------------------------

# Row 31 - Hate Crimes: C31 goes from 1 -> "0" (same text/style as D31)
$ws.Range("D31").Copy($ws.Range("C31"))

# Row 33 - Traffic Fatalities: D33 goes from 1 -> "0" (same as C33)
# and E33 goes from -100 -> "***.*" (same as E31)
$ws.Range("C33").Copy($ws.Range("D33"))
$ws.Range("E31").Copy($ws.Range("E33"))
